$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: insert a new price block (3 quality rows) at the top of the
# data for this market/product, pushing every existing row down by 3.
# Old rows 744-840 become rows 747-843; dimension grows from T840 to T843.
$ws.Rows.Item(744).Resize(3).Insert()

$unidad = "`$/bandeja 4 kilos"

# Shared column values for the 3 inserted rows (only L/M/N/O/P/R/S/D-quality vary).
$common = @{
    A = 8
    B = "Terminal La Palmera de La Serena"
    C = "Coquimbo"
    D = 44776
    E = 4
    F = "Fruta"
    G = 100108
    H = "Tropicales y subtropicales"
    I = 100108002
    J = "Mango"
    K = "Sin especificar"
    N = 10500
    O = 11000
    P = 10750
    Q = $unidad
    R = "Brasil"
    S = 2688
    T = 4
}

$grades = @(
    @{ Row = 744; L = "Especial"; M = 512 },
    @{ Row = 745; L = "Primera";  M = 500 },
    @{ Row = 746; L = "Segunda";  M = 512 }
)

foreach ($grade in $grades) {
    $r = $grade.Row
    $ws.Cells.Item($r, 1).Value = $common.A
    $ws.Cells.Item($r, 2).Value = $common.B
    $ws.Cells.Item($r, 3).Value = $common.C
    $ws.Cells.Item($r, 4).Value = $common.D
    $ws.Cells.Item($r, 5).Value = $common.E
    $ws.Cells.Item($r, 6).Value = $common.F
    $ws.Cells.Item($r, 7).Value = $common.G
    $ws.Cells.Item($r, 8).Value = $common.H
    $ws.Cells.Item($r, 9).Value = $common.I
    $ws.Cells.Item($r, 10).Value = $common.J
    $ws.Cells.Item($r, 11).Value = $common.K
    $ws.Cells.Item($r, 12).Value = $grade.L
    $ws.Cells.Item($r, 13).Value = $grade.M
    $ws.Cells.Item($r, 14).Value = $common.N
    $ws.Cells.Item($r, 15).Value = $common.O
    $ws.Cells.Item($r, 16).Value = $common.P
    $ws.Cells.Item($r, 17).Value = $common.Q
    $ws.Cells.Item($r, 18).Value = $common.R
    $ws.Cells.Item($r, 19).Value = $common.S
    $ws.Cells.Item($r, 20).Value = $common.T
}
